# Update the VIN column (A2:A5) with the new shared VIN value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A5").Value = "DDDKN3DD&E"

# Reflect the new selection captured in the saved file (B9:B10, active cell B10).
$ws.Range("B9:B10").Select()
